$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The pretty-printed (json.dumps(questions, indent=4)) replacement text that
# now lives in A1 (it used to live in A2 as a compact Python dict repr).
$questionsText = @'
questions = [
    {
        "title": "Which of the following modules contains architectures such as the ResNet50 and VGG19?",
        "ques_type": 2,
        "options": [
            "tensorflow.keras.applications",
            "tensorflow.keras.preprocessing",
            "tensorflow.keras.models",
            "tensorflow.keras.layers"
        ],
        "score": "tensorflow.keras.applications"
    },
    {
        "title": "Given the sequential model summary for a Keras sequential model below, what is the code for the third layer?      Model: \"model\"   _________________________________________________________________   Layer (type)                Output Shape             Param #     =================================================================   input_1 (InputLayer)        [(None, None)]           0           _________________________________________________________________   embedding (Embedding)       (None, None, 128)        2560000     _________________________________________________________________   bidirectional (Bidirectional (None, None, 128)        98816       _________________________________________________________________   bidirectional_1 (Bidirectional (None, 128)              98816       _________________________________________________________________   dense (Dense)               (None, 1)                129         =================================================================   Total params: 2,757,761   Trainable params: 2,757,761   Non-trainable params: 0   _________________________________________________________________",
        "ques_type": 2,
        "options": [
            "layers.Bidirectional(layers.LSTM(128, return_sequences=True))",
            "layers.Bidirectional(layers.LSTM(32))",
            "layers.Bidirectional(layers.LSTM(128))",
            "layers.Bidirectional(layers.LSTM(64, return_sequences=True))"
        ],
        "score": "layers.Bidirectional(layers.LSTM(64, return_sequences=True))"
    },
    {
        "title": "The target variable of a dataset with 1024 samples and 6 classes has the shape (1024, 6).  What loss function must be used when compiling the model?",
        "ques_type": 2,
        "options": [
            "mean-squared-error",
            "sparse-categorical-crossentropy",
            "categorical-crossentropy",
            "mean-absolute-error"
        ],
        "score": "categorical-crossentropy"
    },
    {
        "title": "Your dataset has images of the size 32x40 with 2 classes. When you run the code below, an error pops up saying \u201cExpected shape (None, 32).\u201d Why is this error occurring? model = tf.keras.Sequential([          tf.keras.layers.Flatten(input_shape=(32, 40)),          tf.keras.layers.Dense(128, activation='relu'),          tf.keras.layers.Dense(32)           ])",
        "ques_type": 2,
        "options": [
            "The second dense layer must have 32*40 neurons.",
            "The activation function must not be ReLU.",
            "The input shape must be (40, 32).",
            "The last dense layer must have only 2 neurons."
        ],
        "score": "The last dense layer must have only 2 neurons."
    }
]
'@

# A2 used to hold this text as a shared string; it is being removed so the
# sheet's used range shrinks back down to just A1.
$ws.Range("A2").ClearContents()

# A1 used to hold a placeholder 0 with a bold/centered/bordered style (s="1"
# in the original sheet). Clear that formatting back to the workbook default
# before writing the new text so the cell carries no style index.
$ws.Range("A1").ClearFormats()
$ws.Range("A1").Value = $questionsText

# The new text contains embedded newlines, which would otherwise make the
# engine auto-expand the row height (customHeight). AutoFit restores the
# default (non-custom) row height so no ht/customHeight attributes linger.
$ws.Rows.Item(1).AutoFit()
